$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.84406969457198744
$ws.Range("T1").Value = 0.95030951142263598
$ws.Range("BM1").Value = 0.79424805676238663
$ws.Range("D2").Value = 0.93595585834606732
$ws.Range("BL2").Value = 0.91419191008359224
$ws.Range("BP2").Value = 0.9805945032983443
$ws.Range("B3").Value = 0.97363867757716016
$ws.Range("U3").Value = 0.7638071239718851
$ws.Range("BF3").Value = 0.91260585487920265
$ws.Range("C4").Value = 0.78997406186949171
$ws.Range("E4").Value = 0.98473751954392119
$ws.Range("C5").Value = 0.8141080313907787
$ws.Range("F5").Value = 0.92374416831700756
$ws.Range("AP5").Value = 0.91094817366262215
$ws.Range("D6").Value = 0.80232046515208633
$ws.Range("F7").Value = 0.74942714010387768
$ws.Range("I7").Value = 0.96176481605732678
$ws.Range("F8").Value = 0.80215204443526744
$ws.Range("G8").Value = 0.88293810122818506
$ws.Range("J8").Value = 0.77949204753510437
$ws.Range("S8").Value = 0.6403740112219718
$ws.Range("K9").Value = 0.64939960975492173
$ws.Range("X9").Value = 0.75395245632674923
$ws.Range("K10").Value = 0.68291851766607237
$ws.Range("L10").Value = 0.64842767621753361
$ws.Range("M11").Value = 0.94444497242210668
$ws.Range("AM11").Value = 0.98485603880536399
$ws.Range("M12").Value = 0.77484136047155971
$ws.Range("Y12").Value = 0.96051978076466527
$ws.Range("O13").Value = 0.94175426711000543
$ws.Range("M14").Value = 0.81200417258731983
$ws.Range("P14").Value = 0.98304423532742646
$ws.Range("Q15").Value = 0.67937858356756475
$ws.Range("Q18").Value = 0.70126401934584115
$ws.Range("AY18").Value = 0.9620938048640576
$ws.Range("N19").Value = 0.99112050760299464
$ws.Range("Q19").Value = 0.97485589968068642
$ws.Range("U19").Value = 0.69731825001117009
$ws.Range("BG19").Value = 0.82628623597971773
$ws.Range("R20").Value = 0.96693744910405399
$ws.Range("BH20").Value = 0.79147785517744418
$ws.Range("W21").Value = 0.77948378610572167
$ws.Range("W22").Value = 0.85240086646409607
$ws.Range("E23").Value = 0.89104935177213351
$ws.Range("X23").Value = 0.87835381125561773
$ws.Range("Y23").Value = 0.77169794416943005
$ws.Range("V24").Value = 0.93600897065996969
$ws.Range("Y24").Value = 0.77841504908580617
$ws.Range("Z24").Value = 0.65645607457144606
$ws.Range("AK25").Value = 0.88417838699252438
$ws.Range("BH25").Value = 0.98050585764214748
$ws.Range("AB26").Value = 0.82389633411557117
$ws.Range("Z27").Value = 0.76660418488388549
$ws.Range("AD28").Value = 0.95991987228892195
$ws.Range("AA29").Value = 0.96886329540628247
$ws.Range("N30").Value = 0.84338280305506719
$ws.Range("O30").Value = 0.73308411711287214
$ws.Range("Z30").Value = 0.83690183053056255
$ws.Range("AC30").Value = 0.85747038282458321
$ws.Range("AC31").Value = 0.87377571193761683
$ws.Range("AG31").Value = 0.82541432547104177
$ws.Range("AH32").Value = 0.94902016754040086
$ws.Range("AT32").Value = 0.8007400252354917
$ws.Range("P33").Value = 0.69228380132272527
$ws.Range("AG34").Value = 0.95359307852569664
$ws.Range("AI34").Value = 0.77053135277114881
$ws.Range("AJ34").Value = 0.62671915148163171
$ws.Range("AG35").Value = 0.74165524745943423
$ws.Range("AJ35").Value = 0.83950283467304654
$ws.Range("AA36").Value = 0.76559541934169184
$ws.Range("T37").Value = 0.67178769225500434
$ws.Range("AI37").Value = 0.97932967987832953
$ws.Range("AJ37").Value = 0.68099622426497941
$ws.Range("AL37").Value = 0.85992271874923687
$ws.Range("AJ38").Value = 0.67089198889248802
$ws.Range("AN38").Value = 0.72258194447667756
$ws.Range("AR38").Value = 0.98378272731212535
$ws.Range("AK39").Value = 0.94006255623494117
$ws.Range("AO39").Value = 0.97380736535945567
$ws.Range("BP39").Value = 0.74817248509403678
$ws.Range("AM40").Value = 0.86700411280557232
$ws.Range("AP40").Value = 0.71086780245212067
$ws.Range("G42").Value = 0.97992666034121911
$ws.Range("R42").Value = 0.92543101909177072
$ws.Range("AO42").Value = 0.94032464047300601
$ws.Range("AO43").Value = 0.95813480018044284
$ws.Range("AS43").Value = 0.64385653895181605
$ws.Range("BG43").Value = 0.95884419148976407
$ws.Range("E44").Value = 0.98661554830891285
$ws.Range("AB44").Value = 0.61796044767477498
$ws.Range("AE44").Value = 0.76084779889036702
$ws.Range("AS44").Value = 0.87620743180455252
$ws.Range("AU46").Value = 0.96215535586449019
$ws.Range("AV46").Value = 0.67806953366543021
$ws.Range("AS47").Value = 0.82671690517285401
$ws.Range("AV47").Value = 0.74413726811646752
$ws.Range("AX48").Value = 0.69187949848405228
$ws.Range("W49").Value = 0.83899304605799263
$ws.Range("AK49").Value = 0.72479110223193111
$ws.Range("AA50").Value = 0.95731020571439407
$ws.Range("AS51").Value = 0.78136923471395647
$ws.Range("AZ51").Value = 0.84000749833154686
$ws.Range("M52").Value = 0.77012332752753521
$ws.Range("Y52").Value = 0.91985204254382369
$ws.Range("P53").Value = 0.55921137224037021
$ws.Range("BC53").Value = 0.91840703436810966
$ws.Range("E54").Value = 0.95907326201324949
$ws.Range("BA54").Value = 0.58257268606411527
$ws.Range("AU55").Value = 0.91779915081840058
$ws.Range("BD55").Value = 0.9525563666841852
$ws.Range("BE55").Value = 0.9733556204013798
$ws.Range("BH55").Value = 0.85295154108575377
$ws.Range("AV56").Value = 0.93826626379340317
$ws.Range("BB56").Value = 0.95202085132635472
$ws.Range("BE56").Value = 0.95850901312995762
$ws.Range("BK56").Value = 0.91432278171711623
$ws.Range("AX58").Value = 0.93857685347310049
$ws.Range("BE58").Value = 0.79805609938247957
$ws.Range("BE59").Value = 0.88817344948464805
$ws.Range("BF59").Value = 0.54366460010861195
$ws.Range("BH59").Value = 0.80330305898358878
$ws.Range("BL59").Value = 0.69384436723958232
$ws.Range("AN60").Value = 0.88562167873052333
$ws.Range("BG61").Value = 0.57691501486475039
$ws.Range("BI62").Value = 0.78172407869444804
$ws.Range("BI63").Value = 0.80547499391206967
$ws.Range("BJ63").Value = 0.71160179353675157
$ws.Range("BJ64").Value = 0.81607716447014633
$ws.Range("BK64").Value = 0.81454424609749521
$ws.Range("BK65").Value = 0.63099409450510524
$ws.Range("BM66").Value = 0.68992769159465883
$ws.Range("BO66").Value = 0.61217754208275932
$ws.Range("A67").Value = 0.79410381005545827
$ws.Range("BM67").Value = 0.73337844376792249
$ws.Range("BN68").Value = 0.6559943835751183
$ws.Range("BO68").Value = 0.9551444920833021
